$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated "Sample ID" query for the SamplesTab row (B3): the Tumor and
# Analyte Type columns were dropped from the SELECT list.
$newSampleIdQuery = "SELECT`n" +
"    DISTINCT (smp.sample_id) AS ""Sample ID"",`n" +
"    sp.participant_id AS ""Participant ID"", `n" +
"    s.study_name AS ""Study Name"",`n" +
"    s.phs_accession AS Accession`n" +
"FROM `n" +
"    df_participant sp`n" +
"JOIN `n" +
"    df_study s ON sp.""study.phs_accession"" = s.phs_accession`n" +
"JOIN `n" +
"    df_sample smp ON smp.""participant.study_participant_id"" = sp.study_participant_id`n" +
"JOIN`n" +
"    df_diagnosis d ON d.""participant.study_participant_id"" = sp.study_participant_id`n" +
"JOIN`n" +
"    df_program p ON p.program_acronym = s.""program.program_acronym""`n" +
"JOIN`n" +
"    df_file f1 ON f1.""sample.sample_id"" = smp.sample_id`n" +
"JOIN`n" +
"    df_genomic_info gi ON gi.""file.file_id"" = f1.file_id`n" +
"WHERE `n" +
"   s.phs_accession = 'phs001437' AND gi.instrument_model = 'Illumina NextSeq500'`n" +
"ORDER BY `n" +
"    smp.sample_id ASC`n" +
"LIMIT 100;"

$ws.Range("B3").Value = $newSampleIdQuery

# Move the view / selection like it ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B3").Select()
